# Apply the NPOC readme edits:
#  - Add two new header columns: E1 "Sample wt", F1 "Total vol: "
#  - Rename the raw sample names in column A (rows that still hold the
#    original, un-blank-corrected 20220808 run labels) by prefixing them
#    with "TMP_" ahead of re-processing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns
$ws.Range("E1").Value = "Sample wt"
$ws.Range("F1").Value = "Total vol: "

# Rows whose Sample Name (column A) needs the TMP_ prefix added
$rowsToPrefix = @(9,10,11,12,13,14,15,16,17,18,21,22,23,24,25,26,27,28,29,30,33,34,35,36,37,38,39,40,41,42,45,46)

foreach ($r in $rowsToPrefix) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "TMP_" + $cell.Value2
}

# Move the active selection to A47 (matches the saved view state)
$ws.Range("A47").Select()
